$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 95

$ws.Cells.Item($row, 1).Value = "Globo"
$ws.Cells.Item($row, 2).Value = "RJ TV 1"
$ws.Cells.Item($row, 3).Value = "Educação"
$ws.Cells.Item($row, 4).Value = "2025-04-09T12:48"
$ws.Cells.Item($row, 5).Value = "Neutro"
$ws.Cells.Item($row, 6).Value = "Tentativa de homicídio em Campos. Ex-agente penal baleado na cabeça segue internado em estado grave. Crime foi ontem. Após ser atingido, colidiu com muro da creche. Repórter *ao vivo* em frente à delegacia de Guarus. *Havia circulado informação de que as aulas na Creche Francisco Alves Dias haviam sido suspensas, mas essa informação não procede*. Segundo a Secretaria de Educação, os alunos foram liberados mais cedo após a colisão do veículo com o muro*. Setor de Infraestrutura da prefeitura aguarda Polícia Civil liberar o local para os devidos reparos. "
